# Weekly update: insert 3 new price rows (Especial / Primera / Segunda, Perú
# origin, week of 2023-08-?? / serial 45147) at the top of the "Palta" price
# history table, pushing the existing rows down by 3 (1008:1059 -> 1011:1062).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows right before the current row 1008; this shifts the
# existing rows 1008:1059 down to 1011:1062 and copies formatting (incl. the
# date style on column D) down from the row above, exactly like Excel's
# native "Insert Copied/Entire Rows" behaviour.
$ws.Rows("1008:1010").Insert()

# Columns A-C, E-K are identical for every record in this block (same
# market/product/category), so populate them for all three new rows.
for ($r = 1008; $r -le 1010; $r++) {
    $ws.Cells.Item($r, 1).Value  = 7
    $ws.Cells.Item($r, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
    $ws.Cells.Item($r, 3).Value  = "Ñuble"
    $ws.Cells.Item($r, 5).Value  = 16
    $ws.Cells.Item($r, 6).Value  = "Fruta"
    $ws.Cells.Item($r, 7).Value  = 100106
    $ws.Cells.Item($r, 8).Value  = "Oleaginosos"
    $ws.Cells.Item($r, 9).Value  = 100106002
    $ws.Cells.Item($r, 10).Value = "Palta"
    $ws.Cells.Item($r, 11).Value = "Hass"
    $ws.Cells.Item($r, 17).Value = "$/bandeja 10 kilos"
    $ws.Cells.Item($r, 18).Value = "Perú"
    $ws.Cells.Item($r, 20).Value = 10
}

# Row 1008: Especial
$ws.Cells.Item(1008, 4).Value  = 45147
$ws.Cells.Item(1008, 12).Value = "Especial"
$ws.Cells.Item(1008, 13).Value = 60
$ws.Cells.Item(1008, 14).Value = 28000
$ws.Cells.Item(1008, 15).Value = 28000
$ws.Cells.Item(1008, 16).Value = 28000
$ws.Cells.Item(1008, 19).Value = 2800

# Row 1009: Primera
$ws.Cells.Item(1009, 4).Value  = 45147
$ws.Cells.Item(1009, 12).Value = "Primera"
$ws.Cells.Item(1009, 13).Value = 120
$ws.Cells.Item(1009, 14).Value = 24000
$ws.Cells.Item(1009, 15).Value = 25000
$ws.Cells.Item(1009, 16).Value = 24500
$ws.Cells.Item(1009, 19).Value = 2450

# Row 1010: Segunda
$ws.Cells.Item(1010, 4).Value  = 45147
$ws.Cells.Item(1010, 12).Value = "Segunda"
$ws.Cells.Item(1010, 13).Value = 120
$ws.Cells.Item(1010, 14).Value = 21000
$ws.Cells.Item(1010, 15).Value = 22000
$ws.Cells.Item(1010, 16).Value = 21500
$ws.Cells.Item(1010, 19).Value = 2150
